$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Retainers")

# 1. Fix combined grade values in column B for existing rows that now represent a single grade
$ws.Cells.Item(9, 2).Value() = 'S'
$ws.Cells.Item(12, 2).Value() = 'S'
$ws.Cells.Item(52, 2).Value() = 'SR'
$ws.Cells.Item(73, 2).Value() = 'SR'
$ws.Cells.Item(82, 2).Value() = 'SR'

# 2. Append new rows 90-96 that were split out of the combined-grade rows above

# Row 90: Wu Shiren (SSR)
$ws.Cells.Item(89, 1).Copy($ws.Cells.Item(90, 1)) | Out-Null
$ws.Cells.Item(90, 1).Value() = 88
$ws.Cells.Item(90, 2).Value() = 'SSR'
$ws.Cells.Item(90, 3).Value() = 'Wu Shiren'
$ws.Cells.Item(90, 4).Value() = 'Rover'
$ws.Cells.Item(90, 5).Value() = 'Status Level 9'
$ws.Cells.Item(90, 6).Value() = 36
$ws.Cells.Item(90, 7).Value() = 'Ordinary Aptitude: Aptitude +1'
$ws.Cells.Item(90, 8).Value() = 'Art of Chivalry: Aptitude +1'
$ws.Cells.Item(90, 9).Value() = 'Retainer Deployment: When deployed, building earnings +30%'
$ws.Cells.Item(90, 10).Value() = 'Business: Aptitude +1'

# Row 91: Cook Zheng (SSR)
$ws.Cells.Item(89, 1).Copy($ws.Cells.Item(91, 1)) | Out-Null
$ws.Cells.Item(91, 1).Value() = 89
$ws.Cells.Item(91, 2).Value() = 'SSR'
$ws.Cells.Item(91, 3).Value() = 'Cook Zheng'
$ws.Cells.Item(91, 4).Value() = 'Artisan'
$ws.Cells.Item(91, 5).Value() = 'Status Level 13'
$ws.Cells.Item(91, 6).Value() = 36
$ws.Cells.Item(91, 7).Value() = 'Ordinary Aptitude: Aptitude +1'
$ws.Cells.Item(91, 8).Value() = 'Art of Ingenuity: Aptitude +1'
$ws.Cells.Item(91, 9).Value() = 'Retainer Deployment: When deployed, building earnings +30%'
$ws.Cells.Item(91, 10).Value() = 'Business: Aptitude +1'

# Row 92: Calabash Immortal (SSR)
$ws.Cells.Item(89, 1).Copy($ws.Cells.Item(92, 1)) | Out-Null
$ws.Cells.Item(92, 1).Value() = 90
$ws.Cells.Item(92, 2).Value() = 'SSR'
$ws.Cells.Item(92, 3).Value() = 'Calabash Immortal'
$ws.Cells.Item(92, 4).Value() = 'Merchant'
$ws.Cells.Item(92, 5).Value() = 'Bandits Den'
$ws.Cells.Item(92, 6).Value() = 55
$ws.Cells.Item(92, 7).Value() = 'Good Aptitude: Aptitude +2'
$ws.Cells.Item(92, 8).Value() = 'Art of Business: Aptitude +1'
$ws.Cells.Item(92, 9).Value() = 'Retainer Deployment: When deployed, building earnings +30%'
$ws.Cells.Item(92, 10).Value() = 'Business: Aptitude +1'
$ws.Cells.Item(92, 19).Value() = 'Big Profit: Earnings of Merchant Retainer +5% in Trade War and Saltern battles'

# Row 93: Calabash Immortal (UR)
$ws.Cells.Item(89, 1).Copy($ws.Cells.Item(93, 1)) | Out-Null
$ws.Cells.Item(93, 1).Value() = 91
$ws.Cells.Item(93, 2).Value() = 'UR'
$ws.Cells.Item(93, 3).Value() = 'Calabash Immortal'
$ws.Cells.Item(93, 4).Value() = 'Merchant'
$ws.Cells.Item(93, 5).Value() = 'Bandits Den'
$ws.Cells.Item(93, 6).Value() = 55
$ws.Cells.Item(93, 7).Value() = 'Good Aptitude: Aptitude +2'
$ws.Cells.Item(93, 8).Value() = 'Art of Business: Aptitude +1'
$ws.Cells.Item(93, 9).Value() = 'Retainer Deployment: When deployed, building earnings +30%'
$ws.Cells.Item(93, 10).Value() = 'Business: Aptitude +1'
$ws.Cells.Item(93, 19).Value() = 'Big Profit: Earnings of Merchant Retainer +5% in Trade War and Saltern battles'

# Row 94: Xie Lu (SSR)
$ws.Cells.Item(89, 1).Copy($ws.Cells.Item(94, 1)) | Out-Null
$ws.Cells.Item(94, 1).Value() = 92
$ws.Cells.Item(94, 2).Value() = 'SSR'
$ws.Cells.Item(94, 3).Value() = 'Xie Lu'
$ws.Cells.Item(94, 4).Value() = 'Scholar'
$ws.Cells.Item(94, 5).Value() = 'Chef Contest'
$ws.Cells.Item(94, 6).Value() = 55
$ws.Cells.Item(94, 7).Value() = 'Good Aptitude: Aptitude +2'
$ws.Cells.Item(94, 8).Value() = 'Art of Scholarship: Aptitude +1'
$ws.Cells.Item(94, 9).Value() = 'Retainer Deployment: When deployed, building earnings +30%'
$ws.Cells.Item(94, 10).Value() = 'Business: Aptitude +1'
$ws.Cells.Item(94, 19).Value() = 'Big Profit: Earnings of Scholar Retainer +5% in Trade War and Saltern battles'
$ws.Cells.Item(94, 23).Value() = 'Ace Cook: Earnings +5% (For each ace cook recruited, earnings +5%)'
$ws.Cells.Item(94, 24).Value() = 'True Cooking Master Boy: Earnings +5%'

# Row 95: Tang San (SSR)
$ws.Cells.Item(89, 1).Copy($ws.Cells.Item(95, 1)) | Out-Null
$ws.Cells.Item(95, 1).Value() = 93
$ws.Cells.Item(95, 2).Value() = 'SSR'
$ws.Cells.Item(95, 3).Value() = 'Tang San'
$ws.Cells.Item(95, 4).Value() = 'Peasant'
$ws.Cells.Item(95, 5).Value() = 'Soul Land'
$ws.Cells.Item(95, 6).Value() = 55
$ws.Cells.Item(95, 7).Value() = 'Good Aptitude: Aptitude +2'
$ws.Cells.Item(95, 8).Value() = 'Art of Farming: Aptitude +1'
$ws.Cells.Item(95, 9).Value() = 'Retainer Deployment: When deployed, building earnings +30%'
$ws.Cells.Item(95, 10).Value() = 'Business: Aptitude +1'
$ws.Cells.Item(95, 11).Value() = 'Resilient Grass: Aptitude +3'
$ws.Cells.Item(95, 12).Value() = 'Spider Lances: Aptitude +3'
$ws.Cells.Item(95, 13).Value() = 'Incarnation of Haitian: Aptitude +4'
$ws.Cells.Item(95, 14).Value() = 'Ult Art of Farming: Aptitude +3'
$ws.Cells.Item(95, 19).Value() = 'Prosperity: Earnings of Peasant Retainer +20% in Trade War and Saltern battles'
$ws.Cells.Item(95, 20).Value() = 'Fortunate Farmer: Earnings of all peasant retainers +5%'
$ws.Cells.Item(95, 21).Value() = 'Master of Tang Sect: Beauty Charm +10% when deployed to Bazaar Palace'

# Row 96: Tang San (UR)
$ws.Cells.Item(89, 1).Copy($ws.Cells.Item(96, 1)) | Out-Null
$ws.Cells.Item(96, 1).Value() = 94
$ws.Cells.Item(96, 2).Value() = 'UR'
$ws.Cells.Item(96, 3).Value() = 'Tang San'
$ws.Cells.Item(96, 4).Value() = 'Peasant'
$ws.Cells.Item(96, 5).Value() = 'Soul Land'
$ws.Cells.Item(96, 6).Value() = 55
$ws.Cells.Item(96, 7).Value() = 'Good Aptitude: Aptitude +2'
$ws.Cells.Item(96, 8).Value() = 'Art of Farming: Aptitude +1'
$ws.Cells.Item(96, 9).Value() = 'Retainer Deployment: When deployed, building earnings +30%'
$ws.Cells.Item(96, 10).Value() = 'Business: Aptitude +1'
$ws.Cells.Item(96, 11).Value() = 'Resilient Grass: Aptitude +3'
$ws.Cells.Item(96, 12).Value() = 'Spider Lances: Aptitude +3'
$ws.Cells.Item(96, 13).Value() = 'Incarnation of Haitian: Aptitude +4'
$ws.Cells.Item(96, 14).Value() = 'Ult Art of Farming: Aptitude +3'
$ws.Cells.Item(96, 19).Value() = 'Prosperity: Earnings of Peasant Retainer +20% in Trade War and Saltern battles'
$ws.Cells.Item(96, 20).Value() = 'Fortunate Farmer: Earnings of all peasant retainers +5%'
$ws.Cells.Item(96, 21).Value() = 'Master of Tang Sect: Beauty Charm +10% when deployed to Bazaar Palace'
